$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 862
$ws.Range("J2").Value = 294.14285
$ws.Range("L2").Value = 294.14285
$ws.Range("N2").Value = -520.14285

$ws.Range("H33").Value = 525.75
$ws.Range("I33").Value = 401.08334
$ws.Range("J33").Value = 899.75
$ws.Range("K33").Value = 401.08334
$ws.Range("L33").Value = 899.75
$ws.Range("M33").Value = -172.08334
$ws.Range("N33").Value = -1357.75

$ws.Range("H112").Value = 2189.6792
$ws.Range("J112").Value = 2318.8542
$ws.Range("L112").Value = 6956.562600000001
$ws.Range("N112").Value = -9172.562600000001

$ws.Range("H116").Value = 5294652.5
$ws.Range("I116").Value = 9262584
$ws.Range("J116").Value = 4078.2222
$ws.Range("K116").Value = 9262584
$ws.Range("L116").Value = 4078.2222
$ws.Range("M116").Value = -9259142
$ws.Range("N116").Value = -10962.2222

$ws.Range("H120").Value = 60000
$ws.Range("J120").Value = 60000
$ws.Range("L120").Value = 60000
$ws.Range("N120").Value = -69676

$ws.Range("H132").Value = 227237.78
$ws.Range("I132").Value = 260137.84
$ws.Range("J132").Value = 9274.875
$ws.Range("K132").Value = 780413.52
$ws.Range("L132").Value = 27824.625
$ws.Range("M132").Value = -777883.52
$ws.Range("N132").Value = -32884.625

$ws.Range("H137").Value = 7148.923
$ws.Range("I137").Value = 6910.222
$ws.Range("K137").Value = 20730.666
$ws.Range("M137").Value = -18180.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1615721.2
$ws.Range("I32").Value = 3981.3103
$ws.Range("K32").Value = 3981.3103
$ws.Range("M32").Value = -3694.3103

$ws.Range("H45").Value = 1489.4
$ws.Range("I45").Value = 984.8570999999999
$ws.Range("K45").Value = 984.8570999999999
$ws.Range("M45").Value = -607.8570999999999

$ws.Range("H63").Value = 5415
$ws.Range("I63").Value = 4664
$ws.Range("J63").Value = 5707.0557
$ws.Range("K63").Value = 4664
$ws.Range("L63").Value = 5707.0557
$ws.Range("M63").Value = -3978
$ws.Range("N63").Value = -7079.0557

$ws.Range("H66").Value = 5415
$ws.Range("I66").Value = 4664
$ws.Range("J66").Value = 5707.0557
$ws.Range("K66").Value = 23320
$ws.Range("L66").Value = 28535.2785
$ws.Range("M66").Value = -19888
$ws.Range("N66").Value = -35399.2785

$ws.Range("H110").Value = 4553.514
$ws.Range("I110").Value = 1851.625
$ws.Range("K110").Value = 1851.625
$ws.Range("M110").Value = 193.375

$ws.Range("H122").Value = 3714.1875
$ws.Range("I122").Value = 3379
$ws.Range("J122").Value = 5166.6665
$ws.Range("K122").Value = 10137
$ws.Range("L122").Value = 15499.9995
$ws.Range("M122").Value = -7687
$ws.Range("N122").Value = -20399.9995

$ws.Range("H132").Value = 557354.3
$ws.Range("I132").Value = 560353.7
$ws.Range("K132").Value = 1681061.1
$ws.Range("M132").Value = -1678531.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 6559.6265
$ws.Range("I99").Value = 6134.1113
$ws.Range("K99").Value = 6134.1113
$ws.Range("M99").Value = -4636.1113

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7007.5713
$ws.Range("I31").Value = 1449.8334
$ws.Range("J31").Value = 11175.875
$ws.Range("K31").Value = 1449.8334
$ws.Range("L31").Value = 11175.875
$ws.Range("M31").Value = -1154.8334
$ws.Range("N31").Value = -11765.875

$ws.Range("H34").Value = 7007.5713
$ws.Range("I34").Value = 1449.8334
$ws.Range("J34").Value = 11175.875
$ws.Range("K34").Value = 1449.8334
$ws.Range("L34").Value = 11175.875
$ws.Range("M34").Value = -1247.8334
$ws.Range("N34").Value = -11579.875

$ws.Range("H58").Value = 34494904
$ws.Range("I58").Value = 41675508
$ws.Range("K58").Value = 41675508
$ws.Range("M58").Value = -41675305

$ws.Range("H92").Value = 18144.875
$ws.Range("J92").Value = 18144.875
$ws.Range("L92").Value = 18144.875
$ws.Range("N92").Value = -23136.875

$ws.Range("H99").Value = 6176507.5
$ws.Range("I99").Value = 12349020
$ws.Range("K99").Value = 12349020
$ws.Range("M99").Value = -12347522

$ws.Range("H126").Value = 6176507.5
$ws.Range("I126").Value = 12349020
$ws.Range("K126").Value = 37047060
$ws.Range("M126").Value = -37044590

$ws.Range("H132").Value = 14632.4
$ws.Range("I132").Value = 5321.923
$ws.Range("J132").Value = 31923.285
$ws.Range("K132").Value = 15965.769
$ws.Range("L132").Value = 95769.855
$ws.Range("M132").Value = -13435.769
$ws.Range("N132").Value = -100829.855

$ws.Range("H134").Value = 50009216
$ws.Range("I134").Value = 62504724
$ws.Range("K134").Value = 187514172
$ws.Range("M134").Value = -187511637

$ws.Range("H136").Value = 34494904
$ws.Range("I136").Value = 41675508
$ws.Range("K136").Value = 125026524
$ws.Range("M136").Value = -125023974

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7983373
$ws.Range("I4").Value = 7983373
$ws.Range("K4").Value = 23950119
$ws.Range("M4").Value = -23950007

$ws.Range("H21").Value = 261
$ws.Range("I21").Value = 541.5
$ws.Range("J21").Value = 120.75
$ws.Range("K21").Value = 1624.5
$ws.Range("L21").Value = 362.25
$ws.Range("M21").Value = -1451.5
$ws.Range("N21").Value = -708.25

$ws.Range("H86").Value = 2566.6667
$ws.Range("I86").Value = 2532.6667
$ws.Range("J86").Value = 2600.6667
$ws.Range("K86").Value = 7598.000100000001
$ws.Range("L86").Value = 7802.000100000001
$ws.Range("M86").Value = -6412.000100000001
$ws.Range("N86").Value = -10174.0001

$ws.Range("H89").Value = 2566.6667
$ws.Range("I89").Value = 2532.6667
$ws.Range("J89").Value = 2600.6667
$ws.Range("K89").Value = 22794.0003
$ws.Range("L89").Value = 23406.0003
$ws.Range("M89").Value = -16866.0003
$ws.Range("N89").Value = -35262.0003

$ws.Range("H122").Value = 5085.609
$ws.Range("J122").Value = 6561.8237
$ws.Range("L122").Value = 59056.4133
$ws.Range("N122").Value = -63956.4133

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 64.28570999999999
$ws.Range("I2").Value = 47.125
$ws.Range("J2").Value = 87.166664
$ws.Range("K2").Value = 47.125
$ws.Range("L2").Value = 87.166664
$ws.Range("M2").Value = 65.875
$ws.Range("N2").Value = -313.166664

$ws.Range("H5").Value = 4879.75
$ws.Range("I5").Value = 4840
$ws.Range("J5").Value = 4999
$ws.Range("K5").Value = 4840
$ws.Range("L5").Value = 4999
$ws.Range("M5").Value = -4728
$ws.Range("N5").Value = -5223

$ws.Range("H57").Value = 19999
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H93").Value = 95194
$ws.Range("J93").Value = 95194
$ws.Range("L93").Value = 95194
$ws.Range("N93").Value = -98938

$ws.Range("H113").Value = 7507.1113
$ws.Range("I113").Value = 4055.4
$ws.Range("J113").Value = 11821.75
$ws.Range("K113").Value = 4055.4
$ws.Range("L113").Value = 11821.75
$ws.Range("M113").Value = -1885.4
$ws.Range("N113").Value = -16161.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10794.45
$ws.Range("I7").Value = 8399.571
$ws.Range("J7").Value = 16382.5
$ws.Range("K7").Value = 8399.571
$ws.Range("L7").Value = 16382.5
$ws.Range("M7").Value = -8287.571
$ws.Range("N7").Value = -16606.5

$ws.Range("H40").Value = 22906.588
$ws.Range("I40").Value = 32675.363
$ws.Range("J40").Value = 4997.1665
$ws.Range("K40").Value = 32675.363
$ws.Range("L40").Value = 4997.1665
$ws.Range("M40").Value = -32539.363
$ws.Range("N40").Value = -5269.1665

$ws.Range("H55").Value = 4059.04
$ws.Range("I55").Value = 2192.4375
$ws.Range("J55").Value = 7377.4443
$ws.Range("K55").Value = 2192.4375
$ws.Range("L55").Value = 7377.4443
$ws.Range("M55").Value = -2019.4375
$ws.Range("N55").Value = -7723.4443

$ws.Range("H61").Value = 4354.325
$ws.Range("I61").Value = 3538.4595
$ws.Range("K61").Value = 3538.4595
$ws.Range("M61").Value = -3336.4595

$ws.Range("H113").Value = 4354.325
$ws.Range("I113").Value = 3538.4595
$ws.Range("K113").Value = 3538.4595
$ws.Range("M113").Value = -1368.4595

$ws.Range("H126").Value = 10794.45
$ws.Range("I126").Value = 8399.571
$ws.Range("J126").Value = 16382.5
$ws.Range("K126").Value = 25198.713
$ws.Range("L126").Value = 49147.5
$ws.Range("M126").Value = -22728.713
$ws.Range("N126").Value = -54087.5

$ws.Range("H132").Value = 2505.7144
$ws.Range("I132").Value = 2505.7144
$ws.Range("K132").Value = 7517.1432
$ws.Range("M132").Value = -4987.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1794.8462
$ws.Range("I107").Value = 2017.1111
$ws.Range("J107").Value = 1294.75
$ws.Range("K107").Value = 6051.3333
$ws.Range("L107").Value = 3884.25
$ws.Range("M107").Value = -4131.3333
$ws.Range("N107").Value = -7724.25

$ws.Range("H126").Value = 3002.842
$ws.Range("I126").Value = 1843.6875
$ws.Range("K126").Value = 5531.0625
$ws.Range("M126").Value = -3061.0625

$ws.Range("H132").Value = 12767.857
$ws.Range("I132").Value = 8103.2104
$ws.Range("K132").Value = 24309.6312
$ws.Range("M132").Value = -21779.6312

$ws.Range("H136").Value = 15637154
$ws.Range("I136").Value = 17865908
$ws.Range("K136").Value = 53597724
$ws.Range("M136").Value = -53595174

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
